# Quarterly report "rolling window" update:
# Drop the oldest quarter (1399/06) and append the newest quarter (1401/12) -
# every quarter column (E:N) shifts one step to the left, and the new
# quarter's figures land in column N. Applies to both header rows (8, 24)
# and every data row (10-20, 26-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header labels (row 8 and row 24) ---------------------------------
$quarters = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

for ($i = 0; $i -lt $quarters.Length; $i++) {
    $col = 5 + $i
    $ws.Cells.Item(8, $col).Value = $quarters[$i]
    $ws.Cells.Item(24, $col).Value = $quarters[$i]
}

# --- Data rows ----------------------------------------------------------
# Each entry: row number -> 10 values for columns E..N (after the shift)
$rowData = @{
    10 = @(0,0,0,0,0,0,0,0,0,0)
    11 = @(0,0,0,0,0,0,0,0,0,0)
    12 = @(0,0,0,0,0,0,0,0,0,0)
    13 = @(0,0,0,0,0,0,0,0,0,0)
    14 = @(0,0,0,0,0,0,0,0,0,0)
    15 = @(845,-845,0,574,-90,-100,157,-155,530,-136)
    16 = @(364,623,477,491,531,540,508,519,476,0)
    17 = @(10405,9915,17644,17235,13794,17655,25887,23782,39822,11598)
    18 = @(0,0,0,0,0,0,0,0,0,0)
    19 = @(15248,16254,14813,22485,18943,24896,19232,42287,12279,55009)
    20 = @(26862,25947,32934,40785,33178,42991,45784,66433,53107,66471)
    26 = @(134,106,119,135,131,109,130,135,137,137)
    27 = @(127,146,129,113,112,137,123,120,116,116)
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 5 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
